$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '65.938.22'
Set-TextValue 'E2' '  -2.45%  '
Set-TextValue 'D3' '3.483.79'
Set-TextValue 'E3' '  +1.17%  '
Set-TextValue 'D5' '582.52'
Set-TextValue 'E5' '  -1.00%  '
Set-TextValue 'D6' '173.08'
Set-TextValue 'E6' '  -3.26%  '
Set-TextValue 'E7' '  +0.05%  '
Set-TextValue 'E8' '  -1.44%  '
Set-TextValue 'D9' '3.481.75'
Set-TextValue 'E9' '  +1.22%  '
Set-TextValue 'E10' '  -5.33%  '
Set-TextValue 'E11' '  -1.59%  '
Set-TextValue 'D12' '0.410'
Set-TextValue 'E12' '  -3.67%  '
Set-TextValue 'D13' '4.086.16'
Set-TextValue 'E13' '  +1.14%  '
Set-TextValue 'E14' '  +1.31%  '
Set-TextValue 'D15' '29.88'
Set-TextValue 'E15' '  -6.38%  '
Set-TextValue 'D16' '66.002.38'
Set-TextValue 'E16' '  -2.30%  '
Set-TextValue 'E17' '  -2.80%  '
Set-TextValue 'D18' '3.486.58'
Set-TextValue 'E18' '  +1.22%  '
Set-TextValue 'D19' '5.93'
Set-TextValue 'E19' '  -3.04%  '
Set-TextValue 'D20' '13.93'
Set-TextValue 'E20' '  -0.11%  '
Set-TextValue 'D21' '366.72'
Set-TextValue 'E21' '  -4.90%  '
Set-TextValue 'E22' '  -1.06%  '
Set-TextValue 'D23' '72.91'
Set-TextValue 'E23' '  +2.26%  '
Set-TextValue 'E24' '  +0.11%  '
Set-TextValue 'E25' '  +6.21%  '
Set-TextValue 'D26' '0.533'
Set-TextValue 'E26' '  +0.58%  '
Set-TextValue 'D27' '9.65'
Set-TextValue 'E27' '  -5.08%  '
Set-TextValue 'D28' '0.179'
Set-TextValue 'E28' '  +2.22%  '
Set-TextValue 'E29' '  +0.06%  '
Set-TextValue 'D30' '24.07'
Set-TextValue 'E30' '  +2.62%  '
Set-TextValue 'E32' '  -2.55%  '
Set-TextValue 'D33' '1.00'
Set-TextValue 'E33' '  +0.05%  '
Set-TextValue 'D34' '7.15'
Set-TextValue 'E34' '  -0.71%  '
Set-TextValue 'E35' '  -5.71%  '
Set-TextValue 'E36' '  -1.00%  '
Set-TextValue 'D37' '29.78'
Set-TextValue 'E37' '  +15.87%  '
Set-TextValue 'D38' '160.89'
Set-TextValue 'E38' '  -0.21%  '
Set-TextValue 'D39' '0.889'
Set-TextValue 'E39' '  +1.13%  '
Set-TextValue 'D40' '2.820.51'
Set-TextValue 'E40' '  +4.80%  '
Set-TextValue 'D41' '1.75'
Set-TextValue 'E41' '  -5.04%  '
Set-TextValue 'D42' '6.50'
Set-TextValue 'E42' '  -1.50%  '
Set-TextValue 'E43' '  -6.68%  '
Set-TextValue 'D44' '4.45'
Set-TextValue 'E44' '  -1.48%  '
Set-TextValue 'D45' '0.0684'
Set-TextValue 'E45' '  -3.56%  '
Set-TextValue 'E46' '  -2.71%  '
Set-TextValue 'D47' '24.12'
Set-TextValue 'E47' '  -6.60%  '
Set-TextValue 'B48' 'Bittensor'
Set-TextValue 'C48' 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue 'D48' '328.71'
Set-TextValue 'E48' '  +1.15%  '
Set-TextValue 'B49' 'VeChain'
Set-TextValue 'C49' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D49' '0.0288'
Set-TextValue 'E49' '  -2.38%  '
Set-TextValue 'E50' '  -2.00%  '
Set-TextValue 'E51' '  -2.45%  '
